{"js": "// Quarterly report table update: bump the \"visit 1\" and \"visit 2\" row\n// counts by one (306 -> 307, and the three \"280\" cells -> 281) while\n// leaving the percentages in parentheses untouched.\nconst replacements = [\n  [\"306 (100.0)\", \"307 (100.0)\"],\n  [\"280 (100.0)\", \"281 (100.0)\"],\n  [\"280 (91.8)\", \"281 (91.8)\"],\n];\n\nfor (const [findText, newText] of replacements) {\n  const results = context.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Quarterly report table update: bump the \"visit 1\" and \"visit 2\" row\n# counts by one (306 -> 307, and the three \"280\" cells -> 281) while\n# leaving the percentages in parentheses untouched.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$wdFindWrapNone = 0\n\n$replacements = @(\n    @(\"306 (100.0)\", \"307 (100.0)\"),\n    @(\"280 (100.0)\", \"281 (100.0)\"),\n    @(\"280 (91.8)\", \"281 (91.8)\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n"}
